{"js": "// Commit 23 on \"problem solving\": extend Problem 3 (Predicting Fingers)\n// with the finger-counting conclusion paragraphs.\n//\n// Target shape (see diff):\n//   ...the pinky 6, ring 7, middle 8, first 9 and thumb 10.<NEW RUN> We need\n//   to figure out what finger she will end up on.\n//   [[empty paragraph]]\n//   So if she counts 1 to 10 she will as has already been stated end up on\n//   her thumb.\n//   [[empty paragraph]]\n//   If she counts 1 to 100 she will end up on her thumb because it\n//   increases by 10<bookmark _GoBack, unchanged>\n//   [[empty paragraph]]  (new, trailing)\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph that ends the \"Predicting Fingers\" narrative \u2014 the\n// one whose text ends with \"...thumb 10.\" \u2014 and the final (bookmarked,\n// empty) paragraph right after it, which in this document is the very last\n// paragraph in the body.\nconst items = paragraphs.items;\nconst targetIndex = items.findIndex((p) =>\n  p.text.indexOf(\"the pinky 6, ring 7, middle 8, first 9 and thumb 10.\") !== -1\n);\nif (targetIndex === -1) {\n  throw new Error(\"Could not find the 'Predicting Fingers' paragraph.\");\n}\nconst targetParagraph = items[targetIndex];\nconst lastParagraph = items[items.length - 1];\n\n// 1) Append a new run to the end of the \"...thumb 10.\" paragraph.\ntargetParagraph\n  .getRange(\"End\")\n  .insertText(\" We need to figure out what finger she will end up on.\", \"End\");\n\n// 2) Insert a new empty paragraph, then the \"So if she counts...\" paragraph,\n//    then another empty paragraph \u2014 all immediately after the target\n//    paragraph (and therefore before the old bookmark paragraph).\ntargetParagraph.insertParagraph(\"\", \"After\");\ntargetParagraph.insertParagraph(\n  \"So if she counts 1 to 10 she will as has already been stated end up on her thumb. \",\n  \"After\"\n);\ntargetParagraph.insertParagraph(\"\", \"After\");\n\n// 3) Append the new run onto the (now-shifted) bookmark paragraph \u2014 it was\n//    originally empty except for the _GoBack bookmark, and the diff puts\n//    the new sentence's run before the bookmark inside that same paragraph.\nlastParagraph\n  .getRange(\"Start\")\n  .insertText(\"If she counts 1 to 100 she will end up on her thumb because it increases by 10\", \"Before\");\n\n// 4) Add a trailing empty paragraph after the (former) bookmark paragraph.\nlastParagraph.insertParagraph(\"\", \"After\");\n\nawait context.sync();\n", "ps1": "# Commit 23 on \"problem solving\": extend Problem 3 (Predicting Fingers)\n# with the finger-counting conclusion paragraphs.\n#\n# Target shape (see diff):\n#   ...the pinky 6, ring 7, middle 8, first 9 and thumb 10.<NEW RUN> We need\n#   to figure out what finger she will end up on.\n#   [[empty paragraph]]\n#   So if she counts 1 to 10 she will as has already been stated end up on\n#   her thumb.\n#   [[empty paragraph]]\n#   If she counts 1 to 100 she will end up on her thumb because it\n#   increases by 10<bookmark _GoBack, unchanged>\n#   [[empty paragraph]]  (new, trailing)\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph that ends the \"Predicting Fingers\" narrative \u2014 the\n# one whose text ends with \"...thumb 10.\" \u2014 by scanning the paragraph\n# collection for the tell-tale phrase.\n$targetIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n  $p = $d.Paragraphs.Item($i)\n  if ($p.Range.Text -like \"*the pinky 6, ring 7, middle 8, first 9 and thumb 10.*\") {\n    $targetIndex = $i\n    break\n  }\n}\nif ($targetIndex -eq -1) {\n  throw \"Could not find the 'Predicting Fingers' paragraph.\"\n}\n\n$target = $d.Paragraphs.Item($targetIndex)\n\n# 1) Append a new run to the end of the \"...thumb 10.\" paragraph.\n$target.Range.InsertAfter(\" We need to figure out what finger she will end up on.\")\n\n# 2) Insert a new empty paragraph right after it (re-fetch by index each\n#    time so we keep inserting relative to the freshly-created paragraph,\n#    since Range objects do not auto-follow InsertParagraphAfter).\n$target.Range.InsertParagraphAfter()\n$emptyIndex = $targetIndex + 1\n\n# 3) Insert another paragraph break after that empty paragraph, then fill\n#    the new one with the \"So if she counts...\" sentence.\n$emptyPara = $d.Paragraphs.Item($emptyIndex)\n$emptyPara.Range.InsertParagraphAfter()\n$soIfIndex = $emptyIndex + 1\n$soIfPara = $d.Paragraphs.Item($soIfIndex)\n$soIfPara.Range.InsertAfter(\"So if she counts 1 to 10 she will as has already been stated end up on her thumb. \")\n\n# 4) Insert a trailing empty paragraph after the \"So if...\" paragraph\n#    (this sits right before the original bookmark paragraph).\n$soIfPara.Range.InsertParagraphAfter()\n\n# 5) The bookmark paragraph (just \"_GoBack\") is still the last paragraph in\n#    the document body; insert the new sentence ahead of the bookmark that\n#    already lives there.\n$bookmarkPara = $d.Paragraphs.Item($d.Paragraphs.Count)\n$bookmarkPara.Range.InsertBefore(\"If she counts 1 to 100 she will end up on her thumb because it increases by 10\")\n\n# 6) Add a trailing empty paragraph after the bookmark paragraph.\n$bookmarkPara.Range.InsertParagraphAfter()\n"}
